$wb = $excel.ActiveWorkbook

# --- Sheet "PERMISOS": update row 2 data, remove row 3, narrow some columns ---
$wsPermisos = $wb.Worksheets.Item("PERMISOS")
$wsPermisos.Range("B2").Value = "CAMILA OLMOS"
$wsPermisos.Range("C2").Value = "45"
$wsPermisos.Rows.Item(3).Delete()
$wsPermisos.Columns.Item(2).ColumnWidth = 14
$wsPermisos.Columns.Item(3).ColumnWidth = 5
$wsPermisos.Columns.Item(7).ColumnWidth = 20

# --- Sheet "MMO": clear all data and columns (leaves an empty sheet) ---
$wsMmo = $wb.Worksheets.Item("MMO")
$wsMmo.Cells.Clear()

# --- Sheet "TEM": update student name/DNI and narrow some columns ---
$wsTem = $wb.Worksheets.Item("TEM")
$wsTem.Range("E2").Value = "CAMILA OLMOS"
$wsTem.Range("F2").Value = "45"
$wsTem.Columns.Item(5).ColumnWidth = 14
$wsTem.Columns.Item(6).ColumnWidth = 7

# --- Add sheet protection to "PERMISOS" (per commit: "se agrega proteccion al archivo de excel") ---
$wsPermisos.Protect("C6D3", $false, $true, $true, $true, $true, $true, $true, $true, $true, $true, $true, $true, $true, $true)
